$wb = $excel.ActiveWorkbook

# 1) Rename "Cross references" -> "Database references" (the ambiguous name
#    is replaced by a clear one, per the commit message).
$refSheet = $wb.Worksheets.Item("Cross references")
$refSheet.Name = "Database references"

# 2) The renamed sheet becomes the active/selected tab (last tab), so the
#    previously active "Rate laws" tab loses its tabSelected flag and
#    workbookView.activeTab moves to this sheet's position.
$refSheet.Activate()

# 3) Cosmetic column-width cleanup on "Reactions": columns G:H no longer
#    carry an explicit width override (<col min="7" max="8".../> removed),
#    while the existing header-cell styling on G1:H1 must be preserved.
$reactions = $wb.Worksheets.Item("Reactions")

# Stash G1:H1's current formatting in scratch cells far outside the sheet's
# used range so it can be restored after the column override is cleared.
$reactions.Range("G1:H1").Copy()
$reactions.Range("M1:N1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Clearing formats on the whole column removes the <col> width/format
# override entirely (columns G:H revert to the sheet's default width).
$reactions.Columns("G:H").ClearFormats()

# Restore the original header-cell formatting onto G1:H1.
$reactions.Range("M1:N1").Copy()
$reactions.Range("G1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Remove the scratch cells so they leave no trace (values or formatting).
$reactions.Range("M1:N1").Clear()

# 4) Iterative-calculation setting enabled with a tighter convergence delta.
$excel.Iteration = $true
$excel.MaxChange = 0.0001
